$wb = $excel.ActiveWorkbook

# Helper: remove a hyperlink anchored at a given A1 address on a worksheet (if any)
function Remove-HyperlinkAt($ws, [string]$address) {
    $hls = @($ws.Hyperlinks)
    foreach ($hl in $hls) {
        if ($hl.Range.Address() -eq $address) {
            $hl.Delete()
        }
    }
}

# --- Overview sheet: status rollup text changes from "Ready for handoff" ---
# --- to "Handoff transform failed" (same shared string reused by zh-cn/de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# --- Per-language sheets: zh-cn and de-de share the same edit pattern ---
foreach ($name in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)

    # Status column: "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff File (C2): the handoff transform failed, so there is no
    # handoff file anymore - remove its hyperlink and clear the cell entirely.
    Remove-HyperlinkAt $ws '$C$2'
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime (D2): reset to the zero-value timestamp.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason (H2): "Include" -> "Ignored"
    $ws.Range("H2").Value = "Ignored"
}
